$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 5.41 = 21559.46 pesos`r`n✅ 21559.46 pesos = 5.39 = 969.61 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 185
$wsTasas.Range("O10").Value = 3988.5
$wsTasas.Range("N12").Value = 3998
$wsTasas.Range("O12").Value = 179.805
